$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B values (only rows 18-85 change relative to the original data)
$newB = @{
    18 = 6;   19 = 6;   20 = 6;   21 = 6;   22 = 11;  23 = 11;  24 = 15;  25 = 11;
    26 = 109; 27 = 120; 28 = 137; 29 = 160; 30 = 590; 31 = 635; 32 = 703; 33 = 764;
    34 = 1279;35 = 1354;36 = 1417;37 = 1464;38 = 1835;39 = 1877;40 = 1924;41 = 1955;
    42 = 2121;43 = 2141;44 = 2158;45 = 2174;46 = 2207;47 = 2214;48 = 2217;49 = 2217;
    50 = 2198;51 = 2195;52 = 2187;53 = 2174;54 = 2005;55 = 1991;56 = 1975;57 = 1956;
    58 = 1775;59 = 1753;60 = 1726;61 = 1700;62 = 1448;63 = 1420;64 = 1383;65 = 1342;
    66 = 927; 67 = 883; 68 = 812; 69 = 766; 70 = 355; 71 = 309; 72 = 256; 73 = 230;
    74 = 34;  75 = 22;  76 = 22;  77 = 21;  78 = 10;  79 = 10;  80 = 10;  81 = 10;
    82 = 2;   83 = 2;   84 = 2;   85 = 2
}

for ($r = 2; $r -le 97; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 28

    if ($newB.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value2 = $newB[$r]
    }
}
